$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("Modelo"), matching the bold/bordered header style
# used by the other header cells (copy format from A1, then set the text).
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update numeric metric values in row 2 (re-run / new predictions)
$ws.Range("B2").Value = 0.2422940601436104
$ws.Range("C2").Value = 0.9952643565012547
$ws.Range("D2").Value = 0.3844238686478025

# Add new model description cell F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(learning_rate=0.1))])"
